# Apply cryptos list update (Mon Jun  3 18:20:02 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.196.01"
$ws.Range("E2").Value = "  +2.20%  "
$ws.Range("D3").Value = "3.777.92"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  +0.44%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "622.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.45%  "
$ws.Range("D7").Value = "3.774.90"
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.519"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.70%  "
$ws.Range("E11").Value = "  +1.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.67"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000247"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.16%  "
$ws.Range("D15").Value = "4.413.85"
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("D16").Value = "3.767.61"
$ws.Range("E16").Value = "  -0.26%  "
$ws.Range("D17").Value = "69.230.10"
$ws.Range("E17").Value = "  +2.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.66"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.20%  "
$ws.Range("E19").Value = "  +1.79%  "
$ws.Range("E20").Value = "  -0.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "467.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.63"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.04%  "
$ws.Range("E23").Value = "  +1.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000149"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.42%  "
$ws.Range("E27").Value = "  +4.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.07%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "3.927.85"
$ws.Range("E30").Value = "  +0.22%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.66"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.64%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.78"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.166"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +17.39%  "
$ws.Range("B37").Value = "RenzoRestakedETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D37").Value = "3.729.78"
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.18%  "
$ws.Range("E39").Value = "  +3.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.39"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.52%  "
$ws.Range("E41").Value = "  +1.01%  "
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("B45").Value = "Arweave"
$ws.Range("C45").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.12%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.298"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "152.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "46.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.90"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.20%  "
$ws.Range("E51").Value = "  +1.38%  "